$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Model Accuracy (-0.9, 0.9, 0.9)"
# Add new columns C:G (Market threshold, Market min, Market max, Recall,
# Precision) and update the Accuracy (%) values in column B.
# ---------------------------------------------------------------------------
$wsAcc = $wb.Worksheets.Item("Model Accuracy (-0.9, 0.9, 0.9)")

# Copy the header style from B1 (bold, centered, bordered) onto the new
# header cells before writing their text, so they share the same style
# index as the existing header (s="1") instead of minting a new one.
$wsAcc.Range("B1").Copy()
$wsAcc.Range("C1:G1").PasteSpecial(-4122)

$wsAcc.Range("C1").Value = "Market threshold"
$wsAcc.Range("D1").Value = "Market min"
$wsAcc.Range("E1").Value = "Market max"
$wsAcc.Range("F1").Value = "Recall"
$wsAcc.Range("G1").Value = "Precision"

$wsAcc.Range("B2").Value = 63.20293398533008
$wsAcc.Range("C2").Value = 0.05450546436368681
$wsAcc.Range("D2").Value = -15.55441
$wsAcc.Range("E2").Value = 15.06418
$wsAcc.Range("F2").Value = 0
$wsAcc.Range("G2").Value = 0

$wsAcc.Range("B3").Value = 38.32518337408312
$wsAcc.Range("C3").Value = 0.009583939973006913
$wsAcc.Range("D3").Value = -19.35264
$wsAcc.Range("E3").Value = 13.70093
$wsAcc.Range("F3").Value = 2.144772117962467
$wsAcc.Range("G3").Value = 26.66666666666667

$wsAcc.Range("B4").Value = 92.66503667481662
$wsAcc.Range("C4").Value = 0.04158117063764853
$wsAcc.Range("D4").Value = -18.75314
$wsAcc.Range("E4").Value = 23.33066
$wsAcc.Range("F4").Value = 0
$wsAcc.Range("G4").Value = 0

$wsAcc.Range("B5").Value = 82.09046454767727
$wsAcc.Range("C5").Value = 0.02983403801513819
$wsAcc.Range("D5").Value = -12.78028
$wsAcc.Range("E5").Value = 12.42348
$wsAcc.Range("F5").Value = 0
$wsAcc.Range("G5").Value = 0

$wsAcc.Range("B6").Value = 95.59902200488997
$wsAcc.Range("C6").Value = 0.08368817696170747
$wsAcc.Range("D6").Value = -16.47904
$wsAcc.Range("E6").Value = 14.94325
$wsAcc.Range("F6").Value = 0
$wsAcc.Range("G6").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Confusion Matrix TOTALENERGIES SE (-0.9, 0.9, 0.9)"
# ---------------------------------------------------------------------------
$wsTE = $wb.Worksheets.Item("Confusion Matrix TOTALENERGIES SE (-0.9, 0.9, 0.9)")
$wsTE.Range("B3").Value = 9
$wsTE.Range("C3").Value = 1033
$wsTE.Range("D3").Value = 9

# ---------------------------------------------------------------------------
# Sheet 3: "Confusion Matrix FMC CORP (-0.9, 0.9, 0.9)"
# ---------------------------------------------------------------------------
$wsFMC = $wb.Worksheets.Item("Confusion Matrix FMC CORP (-0.9, 0.9, 0.9)")
$wsFMC.Range("B2").Value = 8
$wsFMC.Range("C2").Value = 17
$wsFMC.Range("D2").Value = 5

$wsFMC.Range("B3").Value = 339
$wsFMC.Range("C3").Value = 592
$wsFMC.Range("D3").Value = 324

$wsFMC.Range("B4").Value = 26
$wsFMC.Range("C4").Value = 43
$wsFMC.Range("D4").Value = 27

# ---------------------------------------------------------------------------
# Sheet 4: "Confusion Matrix BP PLC (-0.9, 0.9, 0.9)"
# ---------------------------------------------------------------------------
$wsBP = $wb.Worksheets.Item("Confusion Matrix BP PLC (-0.9, 0.9, 0.9)")
$wsBP.Range("B3").Value = 40
$wsBP.Range("C3").Value = 1516
$wsBP.Range("D3").Value = 42

# ---------------------------------------------------------------------------
# Sheet 5: "Confusion Matrix STORA ENSO (-0.9, 0.9, 0.9)"
# ---------------------------------------------------------------------------
$wsSE = $wb.Worksheets.Item("Confusion Matrix STORA ENSO (-0.9, 0.9, 0.9)")
$wsSE.Range("B3").Value = 110
$wsSE.Range("C3").Value = 1343
$wsSE.Range("D3").Value = 107

# ---------------------------------------------------------------------------
# Sheet 6: "Confusion Matrix BHP GROUP (-0.9, 0.9, 0.9)"
# ---------------------------------------------------------------------------
$wsBHP = $wb.Worksheets.Item("Confusion Matrix BHP GROUP (-0.9, 0.9, 0.9)")
$wsBHP.Range("B3").Value = 4
$wsBHP.Range("C3").Value = 1564
$wsBHP.Range("D3").Value = 3
